# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header/style pattern used by columns B..H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, matching the formatting used by the
# other header cells (e.g. H1): bold font, thin border, centered/top
# aligned.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Numeric data for rows 2-21
$data = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(9, 9)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(6, 6)
    12 = @(8, 8)
    13 = @(6, 7)
    14 = @(9, 9)
    15 = @(5, 6)
    16 = @(9, 9)
    17 = @(7, 8)
    18 = @(4, 4)
    19 = @(4, 4)
    20 = @(5, 5)
    21 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
